$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 updates
$ws.Range("K3").Value = 1.92
$ws.Range("AR3").Value = 1.93
$ws.Range("AS3").Value = 1.93

# Row 7 updates
$ws.Range("G7").Value = 1.5
$ws.Range("I7").Value = 5.75
$ws.Range("J7").Value = 2.05
$ws.Range("AA7").Value = 8
$ws.Range("AG7").Value = 13
$ws.Range("AH7").Value = 8.5
$ws.Range("AI7").Value = 17
$ws.Range("AM7").Value = 29
